$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('B1').Value = 'قیمت تک :‌   43T'
$ws.Range('C1').Value = 'قیمت دبل :‌   59T'
$ws.Range('B2').Value = 'قیمت تک :‌   48T'
$ws.Range('C2').Value = 'قیمت دبل :‌   64T'
$ws.Range('B3').Value = 'قیمت تک :‌   58T'
$ws.Range('C3').Value = 'قیمت دبل :‌   78T'
$ws.Range('B4').Value = 'قیمت :   75T'
$ws.Range('B5').Value = 'قیمت :   65T'
$ws.Range('B6').Value = 'قیمت گرم :   65T'
$ws.Range('C6').Value = 'قیمت سرد :   65T'
$ws.Range('B7').Value = 'قیمت گرم :   74T'
$ws.Range('B8').Value = 'قیمت گرم :   79T'
$ws.Range('B9').Value = 'قیمت گرم :   89T'
$ws.Range('C9').Value = 'قیمت سرد :   89T'
$ws.Range('B10').Value = 'قیمت گرم :   108T'
$ws.Range('C10').Value = 'قیمت سرد :   108T'
$ws.Range('B11').Value = 'قیمت گرم :   115T'
$ws.Range('B12').Value = 'قیمت گرم :   108T'
$ws.Range('C12').Value = 'قیمت سرد :   108T'
$ws.Range('B13').Value = 'قیمت :   105T'
$ws.Range('B14').Value = 'قیمت :   85T'
$ws.Range('B15').Value = 'قیمت :   69T'
$ws.Range('B16').Value = 'قیمت :   69T'
$ws.Range('B17').Value = 'قیمت :   65T'
$ws.Range('B18').Value = 'قیمت :   98T'
$ws.Range('B19').Value = 'قیمت :   98T'
$ws.Range('B20').Value = 'قیمت :   98T'
$ws.Range('B21').Value = 'قیمت :   115T'
$ws.Range('B22').Value = 'قیمت :   98T'
$ws.Range('B23').Value = 'قیمت :   98T'
$ws.Range('B24').Value = 'قیمت :   79T'
$ws.Range('B25').Value = 'قیمت :   89T'
$ws.Range('B26').Value = 'قیمت :   75T'
$ws.Range('B27').Value = 'قیمت :   85T'
$ws.Range('B28').Value = 'قیمت :   85T'
$ws.Range('B29').Value = 'قیمت :   45T'
$ws.Range('B30').Value = 'قیمت :   55T'
$ws.Range('B31').Value = 'قیمت :   58T'
$ws.Range('B32').Value = 'قیمت :   65T'
$ws.Range('B33').Value = 'قیمت :   110T'
$ws.Range('B34').Value = 'قیمت :   189T'
$ws.Range('B35').Value = 'قیمت :   269T'
$ws.Range('B36').Value = 'قیمت :   68T'
$ws.Range('B37').Value = 'قیمت :   --'
$ws.Range('B38').Value = 'قیمت :   189T'
$ws.Range('B39').Value = 'قیمت :   165T'
$ws.Range('B40').Value = 'قیمت :   170T'
$ws.Range('B41').Value = 'قیمت :   138T'
$ws.Range('B42').Value = 'قیمت :   148T'
$ws.Range('B43').Value = 'قیمت :   165T'
$ws.Range('B44').Value = 'قیمت :   98T'
$ws.Range('B45').Value = 'قیمت :   275T'
$ws.Range('B46').Value = 'قیمت :   89T'
$ws.Range('B47').Value = 'قیمت :   95T'
$ws.Range('B48').Value = 'قیمت :   115T'
$ws.Range('B49').Value = 'قیمت :   118T'
$ws.Range('B50').Value = 'قیمت :   110T'
$ws.Range('B51').Value = 'قیمت :   105T'
$ws.Range('B52').Value = 'قیمت :   120T'
$ws.Range('B53').Value = 'قیمت :   145T'
$ws.Range('B54').Value = 'قیمت :   195T'

$ws.Range('B54').Select()
